$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue "D2" "25.964.11"
Set-TextValue "E2" "  +0.42%  "
Set-TextValue "D3" "1.587.62"
Set-TextValue "E3" "  +0.10%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "210.38"
Set-TextValue "E5" "  +0.21%  "
Set-TextValue "E6" "  -0.16%  "
Set-TextValue "E7" "  +0.04%  "
Set-TextValue "D8" "0.246"
Set-TextValue "E8" "  -0.37%  "
Set-TextValue "D9" "0.0610"
Set-TextValue "E9" "  -1.07%  "
Set-TextValue "D10" "17.92"
Set-TextValue "E10" "  -0.73%  "
Set-TextValue "D11" "0.0808"
Set-TextValue "E11" "  +2.31%  "
Set-TextValue "D12" "1.809.68"
Set-TextValue "E12" "  +0.18%  "
Set-TextValue "D13" "1.585.11"
Set-TextValue "E13" "  -0.04%  "
Set-TextValue "D14" "3.98"
Set-TextValue "E14" "  -1.27%  "
Set-TextValue "D15" "0.510"
Set-TextValue "E15" "  +0.04%  "
Set-TextValue "D16" "25.953.70"
Set-TextValue "E16" "  +0.47%  "
Set-TextValue "D17" "60.00"
Set-TextValue "E17" "  +0.29%  "
Set-TextValue "D18" "0.0₃0720"
Set-TextValue "E18" "  -0.41%  "
Set-TextValue "E19" "  -0.12%  "
Set-TextValue "D20" "198.98"
Set-TextValue "E20" "  +3.85%  "
Set-TextValue "D21" "4.21"
Set-TextValue "E21" "  +0.53%  "
Set-TextValue "D22" "9.17"
Set-TextValue "E22" "  -2.12%  "
Set-TextValue "D23" "5.97"
Set-TextValue "E23" "  +0.55%  "
Set-TextValue "D24" "1.83"
Set-TextValue "E24" "  +7.95%  "
Set-TextValue "D25" "142.57"
Set-TextValue "E25" "  +0.32%  "
Set-TextValue "E26" "  -0.10%  "
Set-TextValue "E27" "  -8.45%  "
Set-TextValue "D28" "15.03"
Set-TextValue "E28" "  -0.43%  "
Set-TextValue "D29" "6.44"
Set-TextValue "E29" "  -0.23%  "
Set-TextValue "E30" "  -0.27%  "
Set-TextValue "E31" "  +0.31%  "
Set-TextValue "D32" "3.11"
Set-TextValue "E32" "  -0.01%  "
Set-TextValue "E33" "  -3.55%  "
Set-TextValue "D34" "1.47"
Set-TextValue "E34" "  -2.08%  "
Set-TextValue "E35" "  +0.10%  "
Set-TextValue "D36" "1.121.94"
Set-TextValue "E36" "  +1.76%  "
Set-TextValue "D37" "0.0163"
Set-TextValue "E37" "  +8.39%  "
Set-TextValue "E38" "  -0.16%  "
Set-TextValue "E39" "  -1.65%  "
Set-TextValue "D40" "0.782"
Set-TextValue "E40" "  +0.49%  "
Set-TextValue "D41" "0.487"
Set-TextValue "E41" "  -3.24%  "
Set-TextValue "D42" "0.776"
Set-TextValue "E42" "  -5.06%  "
Set-TextValue "D43" "1.720.64"
Set-TextValue "E43" "  +0.04%  "
Set-TextValue "D44" "5.09"
Set-TextValue "E44" "  -1.94%  "
Set-TextValue "D45" "91.78"
Set-TextValue "E45" "  -2.25%  "
Set-TextValue "B46" "RenderToken"
Set-TextValue "C46" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D46" "1.48"
Set-TextValue "E46" "  -1.68%  "
Set-TextValue "B47" "Aave"
Set-TextValue "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "53.17"
Set-TextValue "E47" "  -0.10%  "
Set-TextValue "E48" "  -1.07%  "
Set-TextValue "D49" "0.406"
Set-TextValue "E49" "  -0.24%  "
Set-TextValue "E50" "  +0.03%  "
Set-TextValue "D51" "0.0₇0918"
Set-TextValue "E51" "  -17.53%  "
